$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column E (reuses the shared string table, gets next index)
$ws.Range("E1").Value = "Q4に▲100bp"

# Copy the header formatting (style) from D1, which already carries the
# bold/centered/bordered header style used by B1:D1, so no new cell style
# definition is introduced.
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null

# Fill in the new column E values (row 2 - 13)
$ws.Range("E2").Value = 3.78247358508853
$ws.Range("E3").Value = 3.960903665855144
$ws.Range("E4").Value = 4.23709174823011
$ws.Range("E5").Value = 4.195573929455024
$ws.Range("E6").Value = 4.096743681114789
$ws.Range("E7").Value = 4.037086841535631
$ws.Range("E8").Value = 3.971814840580978
$ws.Range("E9").Value = 3.945399758499387
$ws.Range("E10").Value = 3.919457470063083
$ws.Range("E11").Value = 3.906735596429402
$ws.Range("E12").Value = 3.89934144266164
$ws.Range("E13").Value = 3.899053323967651
